# ShapefileAddingIndonesia.xlsx edit
# Commit: "implemented changes to the map window and shapefile to focus on
# the 5 provinces of interest in Indonesia"
#
# The 5 provinces of interest (rows left untouched) are:
#   row 7  -> Central Java
#   row 10 -> East Java
#   row 14 -> Jakarta Special Capital Region
#   row 27 -> South Sulawesi
#   row 30 -> West Java
#
# Every other province (rows 2-6, 8-9, 11-13, 15-26, 28-29, 31-35) has its
# Cases/Recovered/mobility columns (C:J) zeroed out, and the thousands-
# separator style that was applied to some of the Cases/Recovered (C/D)
# cells is cleared at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A is now wide enough to show full province names ------------
$ws.Columns("A").ColumnWidth = 26.666666667

# --- Zero out every province row except the 5 of interest ----------------
$zeroAddr   = "C2:J6,C8:J9,C11:J13,C15:J26,C28:J29,C31:J35"
$zeroAddrCD = "C2:D6,C8:D9,C11:D13,C15:D26,C28:D29,C31:D35"

# Clear the manual "#,##0" style from the Cases/Recovered cells being
# zeroed (matches the author's result where s="1" disappears from those
# cells once the values become 0).
$ws.Range($zeroAddrCD).ClearFormats()

# Set all of Cases, Recovered, Retail, Grocery, Parks, Transit, Workplace,
# Residential to 0 for the non-focus provinces.
$ws.Range($zeroAddr).Value = 0

# --- Update the sheet view / selection to focus near the bottom rows -----
$ws.Range("C31:J35").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
